$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tastes & aromas of varieties")

$ws.Range("B31").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>🍓<strong>Aromas - </strong>Spicy and sweet</p><p>👄<strong>Mouthfeel - </strong>Powerful </p></div>"
$ws.Range("B32").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>🍓<strong>Aromas - </strong>fruity, plum, spicy, cinamon</p><p>👄<strong>Mouthfeel - </strong>Powerful and rich with great fineness</p></div>"
$ws.Range("B33").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>🍓<strong>Aromas - </strong>Fruits and violet</p></div>"
$ws.Range("B35").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red, wine similar to cabernet franc</p><p>👄<strong>Mouthfeel - </strong>Subtil and delicate</p></div>"
$ws.Range("B34").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red, good ageing capacity</p><p>🍓<strong>Aromas - </strong>Spicy and fruity</p></div>"
$ws.Range("B36").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>🍓<strong>Aromas - </strong>Pulpy fruits, rich, sappy with bitter almond notes</p><p>👄<strong>Mouthfeel - </strong>Strong, structured, robust</p></div>"
$ws.Range("B37").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>🍓<strong>Aromas - </strong>Kirsch</p></div>"
$ws.Range("B38").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red, to drink young</p><p>🍓<strong>Aromas - </strong>Fruity</p><p>👄<strong>Mouthfeel - </strong>Fresh</p></div>"
$ws.Range("B39").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red</p><p>👅<strong>Basic tastes - </strong>Dry</p><p>🍓<strong>Aromas - </strong>Floral with cloves notes</p></div>"
$ws.Range("B40").Value2 = "<div><p>🍷<strong>Wine color - </strong>Red, age fast</p><p>👅<strong>Basic tastes - </strong></p><p>🍓<strong>Aromas - </strong>Floral</p><p>👄<strong>Mouthfeel - </strong>Strong</p></div>"
$ws.Range("B41").Value2 = "<div><p>🍷<strong>Wine color - </strong>White</p><p>🍓<strong>Aromas - </strong>Floral, lemongrass, herbs and fruity, apple, grapefruit, with nut taste</p></div>"

$ws.Range("B41").Select()
